$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$tb = $s.Shapes.Item(3)
$tr = $tb.TextFrame.TextRange

$tr.Paragraphs(1).Text = "."
$tr.Paragraphs(1).Text = "SQL Saturday #839"

$tr.Paragraphs(2).Text = "."
$tr.Paragraphs(2).Text = "Virginia Beach, Virginia"

$tr.Paragraphs(3).Text = "."
$tr.Paragraphs(3).Text = "8 June 2019"
